$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# The localization round trip finished: both the zh-cn and de-de handoffs
# came back "in sync with en-US". For each language sheet we:
#   - flip the Status text ("Ready for handoff" -> "Handed back: in sync
#     with en-US") for every row (shared across both sheets/rows),
#   - populate "Latest Target File" (F) / "Latest Handback File" (G) with
#     the returned source + xlf files (same links as the original
#     handoff columns A/D, since the file came back as-is),
#   - stamp "Latest Handback DateTime" (H) with the real handback time
#     instead of the 0001-01-01 placeholder.
# ---------------------------------------------------------------------------

$statusText = "Handed back: in sync with en-US"

$mdUrlBase  = "https://github.com/OpenLocalizationTest/oltest/blob/1fab925f07ed73f77fc0e918e17473b30370ced4/e2e/"
$zhXlfBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3abda281877ffa8cfb8ae9dc1a08c54c5cefc7a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/"
$deXlfBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/30614eb8d6f627ecf3d39e0eff5cd3e87b034f6f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/"

$file1 = "2185879b-a60e-45d6-b123-41b9c4e71e5f"
$file2 = "46977c10-c8ed-4c00-b8b6-df84d21c5099"

$md1 = "$file1.md"
$md2 = "$file2.md"

$zhXlf1 = "$file1.96091fa7731abc6479e34e94c2260793f109df7a.zh-cn.xlf"
$zhXlf2 = "$file2.2b68ecba61e493100692275886bbf9f1e95cd520.zh-cn.xlf"

$deXlf1 = "$file1.96091fa7731abc6479e34e94c2260793f109df7a.de-de.xlf"
$deXlf2 = "$file2.2b68ecba61e493100692275886bbf9f1e95cd520.de-de.xlf"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), ($mdUrlBase + $md1), "", "", $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), ($zhXlfBase + $zhXlf1), "", "", $zhXlf1)

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), ($mdUrlBase + $md2), "", "", $md2)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), ($zhXlfBase + $zhXlf2), "", "", $zhXlf2)

$wsZh.Range("H2").Value = "2016-03-23 10:14:59"
$wsZh.Range("H3").Value = "2016-03-23 10:14:59"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), ($mdUrlBase + $md1), "", "", $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), ($deXlfBase + $deXlf1), "", "", $deXlf1)

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), ($mdUrlBase + $md2), "", "", $md2)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), ($deXlfBase + $deXlf2), "", "", $deXlf2)

$wsDe.Range("H2").Value = "2016-03-23 10:15:09"
$wsDe.Range("H3").Value = "2016-03-23 10:15:09"
